$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.126.06"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").Value = "2.564.02"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'584.15"
$ws.Range("E5").Value = "  +2.64%  "

$ws.Range("D6").Value = "'147.64"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "  +1.38%  "

$ws.Range("E9").Value = "  +2.70%  "

$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("E12").Value = "  +0.78%  "

$ws.Range("D13").Value = "'27.38"
$ws.Range("E13").Value = "  -0.41%  "

$ws.Range("D14").Value = "3.021.78"
$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("D15").Value = "63.040.16"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("D17").Value = "2.542.55"
$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").Value = "'11.35"
$ws.Range("E18").Value = "  -1.27%  "

$ws.Range("D19").Value = "'343.40"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("E21").Value = "  +1.33%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("E23").Value = "  -3.77%  "

$ws.Range("D24").Value = "'66.60"
$ws.Range("E24").Value = "  +1.76%  "

$ws.Range("D25").Value = "2.695.14"
$ws.Range("E25").Value = "  +0.70%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").Value = "'1.63"
$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("D28").Value = "'8.13"
$ws.Range("E28").Value = "  +10.48%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").Value = "'1.48"
$ws.Range("E30").Value = "  -2.11%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.42"
$ws.Range("E31").Value = "  -0.68%  "

$ws.Range("E32").Value = "  +7.07%  "

$ws.Range("D33").Value = "0.0₃0822"
$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("D34").Value = "'460.87"
$ws.Range("E34").Value = "  +11.67%  "

$ws.Range("D35").Value = "'175.67"
$ws.Range("E35").Value = "  -1.59%  "

$ws.Range("E36").Value = "  +2.46%  "

$ws.Range("D37").Value = "'0.407"
$ws.Range("E37").Value = "  +1.97%  "

$ws.Range("E38").Value = "  +0.54%  "

$ws.Range("D39").Value = "'4.52"
$ws.Range("E39").Value = "  +2.72%  "

$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").Value = "'1.75"
$ws.Range("E41").Value = "  -1.15%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").Value = "'150.83"
$ws.Range("E43").Value = "  -1.84%  "

$ws.Range("E44").Value = "  +0.97%  "

$ws.Range("E45").Value = "  -0.80%  "

$ws.Range("D46").Value = "'0.0546"
$ws.Range("E46").Value = "  +4.22%  "

$ws.Range("E47").Value = "  +1.37%  "

$ws.Range("E48").Value = "  +1.12%  "

$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("E50").Value = "  -2.72%  "

$ws.Range("E51").Value = "  +0.34%  "
